$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first three data rows (rows 2-4), shifting remaining data up
$ws.Range("A2:A4").EntireRow.Delete()

# Append new measurement rows at the bottom (rows 19-31)
$ws.Cells.Item(19,1).Value = -2.693654298782349
$ws.Cells.Item(19,2).Value = 12.81020450592041
$ws.Cells.Item(19,3).Value = -6.878878593444824
$ws.Cells.Item(20,1).Value = 4.974119186401367
$ws.Cells.Item(20,2).Value = -6.053498268127441
$ws.Cells.Item(20,3).Value = -4.147007942199707
$ws.Cells.Item(21,1).Value = -2.26168966293335
$ws.Cells.Item(21,2).Value = -17.45715522766113
$ws.Cells.Item(21,3).Value = 0.2626542747020721
$ws.Cells.Item(22,1).Value = -1.576191902160644
$ws.Cells.Item(22,2).Value = 2.079262971878052
$ws.Cells.Item(22,3).Value = 1.673064351081848
$ws.Cells.Item(23,1).Value = -10.3257417678833
$ws.Cells.Item(23,2).Value = 1.127981901168823
$ws.Cells.Item(23,3).Value = 5.69497013092041
$ws.Cells.Item(24,1).Value = 6.999719619750977
$ws.Cells.Item(24,2).Value = -9.156081199645996
$ws.Cells.Item(24,3).Value = -6.504971027374268
$ws.Cells.Item(25,1).Value = 1.670201539993286
$ws.Cells.Item(25,2).Value = -8.627177238464355
$ws.Cells.Item(25,3).Value = 0.0325571447610855
$ws.Cells.Item(26,1).Value = 7.071092128753662
$ws.Cells.Item(26,2).Value = 1.176451444625854
$ws.Cells.Item(26,3).Value = -1.140298962593079
$ws.Cells.Item(27,1).Value = -2.711231231689453
$ws.Cells.Item(27,2).Value = -4.331498622894287
$ws.Cells.Item(27,3).Value = 1.348158717155456
$ws.Cells.Item(28,1).Value = -5.673731803894043
$ws.Cells.Item(28,2).Value = 1.553022384643555
$ws.Cells.Item(28,3).Value = 2.025134325027466
$ws.Cells.Item(29,1).Value = -3.546931266784668
$ws.Cells.Item(29,2).Value = 9.443968772888184
$ws.Cells.Item(29,3).Value = -0.7339006066322327
$ws.Cells.Item(30,1).Value = 0.6049370765686035
$ws.Cells.Item(30,2).Value = 8.170445442199707
$ws.Cells.Item(30,3).Value = -6.977948188781738
$ws.Cells.Item(31,1).Value = 4.74082612991333
$ws.Cells.Item(31,2).Value = -5.404219627380371
$ws.Cells.Item(31,3).Value = -7.638944625854492
